# Updates the "cryptos" price/volume table (columns D and E) with refreshed
# values from the upstream scrape, per the commit
# "Updated cryptos list on Sat Oct  7 22:38:58 UTC 2023 with GitHub Actions".
#
# Column D ("Price") and column E ("Volume(1h)") are stored as plain text
# cells (not numbers) in the source workbook, so each write below forces the
# target cell to Text format before assigning the new string (otherwise
# values like "65.40" or "212.27" would be auto-coerced into numbers and
# lose their trailing zero / text formatting), then restores the cell's
# style to "Normal" so no stray number-format/style is left behind.

function Set-CellText {
    param($Sheet, $Ref, $NewValue)
    $cell = $Sheet.Range($Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "27.958.10"
Set-CellText $ws "E2" "  -0.32%  "
Set-CellText $ws "D3" "1.635.95"
Set-CellText $ws "E3" "  -0.89%  "
Set-CellText $ws "E4" "  -0.06%  "
Set-CellText $ws "D5" "212.27"
Set-CellText $ws "E5" "  -0.80%  "
Set-CellText $ws "E6" "  -1.00%  "
Set-CellText $ws "E7" "  -0.10%  "
Set-CellText $ws "D8" "23.32"
Set-CellText $ws "E8" "  -1.16%  "
Set-CellText $ws "E9" "  -2.75%  "
Set-CellText $ws "E10" "  -0.05%  "
Set-CellText $ws "D11" "0.0884"
Set-CellText $ws "E11" "  +1.06%  "
Set-CellText $ws "D12" "1.868.31"
Set-CellText $ws "E12" "  -0.93%  "
Set-CellText $ws "D13" "1.628.68"
Set-CellText $ws "E13" "  -1.43%  "
Set-CellText $ws "E14" "  -0.68%  "
Set-CellText $ws "E15" "  -0.32%  "
Set-CellText $ws "D16" "65.40"
Set-CellText $ws "E16" "  -0.63%  "
Set-CellText $ws "D17" "27.959.31"
Set-CellText $ws "E17" "  -0.37%  "
Set-CellText $ws "D18" "230.94"
Set-CellText $ws "E18" "  -1.02%  "
Set-CellText $ws "E19" "  -0.15%  "
Set-CellText $ws "D20" "7.55"
Set-CellText $ws "E20" "  -1.66%  "
Set-CellText $ws "E21" "  -0.01%  "
Set-CellText $ws "D22" "4.37"
Set-CellText $ws "E22" "  -0.85%  "
Set-CellText $ws "D23" "10.37"
Set-CellText $ws "E23" "  -3.20%  "
Set-CellText $ws "D24" "2.06"
Set-CellText $ws "E24" "  -3.96%  "
Set-CellText $ws "D25" "154.89"
Set-CellText $ws "E25" "  +1.48%  "
Set-CellText $ws "E26" "  +0.55%  "
Set-CellText $ws "E27" "  -0.85%  "
Set-CellText $ws "E28" "  -1.13%  "
Set-CellText $ws "E29" "  -0.08%  "
Set-CellText $ws "D30" "1.18"
Set-CellText $ws "E30" "  -0.66%  "
Set-CellText $ws "E31" "  -0.46%  "
Set-CellText $ws "D32" "3.40"
Set-CellText $ws "E32" "  +1.68%  "
Set-CellText $ws "D33" "1.407.32"
Set-CellText $ws "E33" "  -2.88%  "
Set-CellText $ws "E34" "  -0.27%  "
Set-CellText $ws "E35" "  -0.25%  "
Set-CellText $ws "E36" "  +9.15%  "
Set-CellText $ws "E37" "  +1.46%  "
Set-CellText $ws "E38" "  +0.43%  "
Set-CellText $ws "D39" "0.564"
Set-CellText $ws "E39" "  +0.72%  "
Set-CellText $ws "E41" "  +0.29%  "
Set-CellText $ws "E42" "  -0.06%  "
Set-CellText $ws "D43" "66.80"
Set-CellText $ws "E43" "  -3.87%  "
Set-CellText $ws "E44" "  +2.11%  "
Set-CellText $ws "E45" "  +0.51%  "
Set-CellText $ws "E46" "  -1.10%  "
Set-CellText $ws "D47" "1.777.69"
Set-CellText $ws "E47" "  -1.02%  "
Set-CellText $ws "D48" "88.01"
Set-CellText $ws "E48" "  -1.30%  "
Set-CellText $ws "E49" "  +0.96%  "
Set-CellText $ws "E50" "  -1.46%  "
Set-CellText $ws "E51" "  -0.42%  "
